# Generate Report for Handoff
# Adds a new tracked file (cdaddbb6-497a-412b-a3ac-2d3837fa9772) as row 7
# to the Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$fileId = "cdaddbb6-497a-412b-a3ac-2d3837fa9772"
$mdName = $fileId + ".md"
$zhXlfName = $fileId + ".6513aaaad878a03b47ecc940d99b792a26b2052a.zh-cn.xlf"
$deXlfName = $fileId + ".6513aaaad878a03b47ecc940d99b792a26b2052a.de-de.xlf"

$dateFmt = "yyyy-mm-dd HH:mm:ss"
# BGR-packed "FF6495ED" (the workbook's existing HyperLink font color) so the
# new hyperlink cells render identically to the pre-existing ones (A2:A6, D2:D6, ...).
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ----------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/b0a3e4ee619ac4914aefb94332c5e488c02f52aa/e2e/" + $mdName, "", "", $mdName)
Style-AsHyperlink $wsOverview.Range("A7")
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-03-25 10:22:09"
$wsOverview.Range("D7").NumberFormat = $dateFmt

# ----------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Reference Tokens |
#   Handoff Reason | Dependency From | Error Detail
# ----------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/fda50076723005c6665817c6bdd7787370fe33f4/e2e/" + $mdName, "", "", $mdName)
Style-AsHyperlink $wsZhCn.Range("A7")
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65d36448d1fcc1af8bd40f4edaf85932c1e2a094/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/" + $zhXlfName, "", "", $zhXlfName)
Style-AsHyperlink $wsZhCn.Range("D7")
$wsZhCn.Range("E7").Value = "2016-03-25 10:22:00"
$wsZhCn.Range("E7").NumberFormat = $dateFmt
$wsZhCn.Range("H7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H7").NumberFormat = $dateFmt
$wsZhCn.Range("J7").Value = "Include"

# ----------------------------------------------------------------
# Sheet "de-de": same column layout as "zh-cn"
# ----------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/b0a3e4ee619ac4914aefb94332c5e488c02f52aa/e2e/" + $mdName, "", "", $mdName)
Style-AsHyperlink $wsDeDe.Range("A7")
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90bf1647a3b94f204f26c5ed4f3ec7187dd8f8b9/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/" + $deXlfName, "", "", $deXlfName)
Style-AsHyperlink $wsDeDe.Range("D7")
$wsDeDe.Range("E7").Value = "2016-03-25 10:22:09"
$wsDeDe.Range("E7").NumberFormat = $dateFmt
$wsDeDe.Range("H7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H7").NumberFormat = $dateFmt
$wsDeDe.Range("J7").Value = "Include"

Write-Output "Row 7 added to Overview, zh-cn, de-de sheets for $fileId"
